# diagnostic.xlsx: write the "disconnected_elements" diagnostic block
# B1 = 0, A2 = 0  (both bold, thin-boxed, centered/top-aligned)
# B2 = "disconnected_elements" (plain, becomes a shared string)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 0
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "disconnected_elements"

$c1 = $ws.Range("B1")
$c1.Font.Bold = $true
$c1.Borders.LineStyle = 1
$c1.VerticalAlignment = -4160
$c1.HorizontalAlignment = -4108

# Clone B1's format onto A2 via copy/paste-special so both cells share the
# exact same style record instead of the engine minting a second, equivalent
# one.
$c1.Copy()
$ws.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false
